$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.87"
$ws.Range("E2").Value = "'-1.13%"
$ws.Range("D3").Value = "'27.50"
$ws.Range("E3").Value = "'5.67%"
$ws.Range("D4").Value = "'5.128"
$ws.Range("E4").Value = "'-1.35%"
$ws.Range("D5").Value = "'0.05689"
$ws.Range("E5").Value = "'1.71%"
$ws.Range("D6").Value = "'6.550"
$ws.Range("E6").Value = "'1.04%"
$ws.Range("D7").Value = "'0.8197"
$ws.Range("E7").Value = "'0.87%"
$ws.Range("D8").Value = "'0.8618"
$ws.Range("E8").Value = "'1.89%"
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").Value = "'0.06948"
$ws.Range("E9").Value = "'0.23%"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "'0.02862"
$ws.Range("E10").Value = "'1.34%"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "'0.09391"
$ws.Range("E11").Value = "'0.11%"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "'0.001510"
$ws.Range("E12").Value = "'0.00%"
$ws.Range("B13").Value = "CoinExToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D13").Value = "'0.04095"
$ws.Range("E13").Value = "'-12.08%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006009"
$ws.Range("E14").Value = "'0.86%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006213"
$ws.Range("E15").Value = "'0.65%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.516"
$ws.Range("E16").Value = "'-2.61%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.010"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.176"
$ws.Range("E18").Value = "'5.90%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3150"
$ws.Range("E19").Value = "'1.22%"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").Value = "'0.1332"
$ws.Range("E20").Value = "'-0.14%"
$ws.Range("D21").Value = "'0.03234"
$ws.Range("E21").Value = "'1.66%"
$ws.Range("D22").Value = "'0.1302"
$ws.Range("E22").Value = "'0.67%"
$ws.Range("D23").Value = "'3.570"
$ws.Range("E23").Value = "'-4.93%"
$ws.Range("D25").Value = "'0.001215"
$ws.Range("E25").Value = "'-2.32%"
$ws.Range("D26").Value = "'0.004467"
$ws.Range("E26").Value = "'-1.79%"
$ws.Range("E27").Value = "'23.02%"
$ws.Range("D28").Value = "'0.0001406"
$ws.Range("E28").Value = "'-27.41%"
$ws.Range("D40").Value = "'0.03715"
$ws.Range("E40").Value = "'1.79%"
$ws.Range("D41").Value = "'0.005921"
$ws.Range("E41").Value = "'-3.50%"
$ws.Range("D42").Value = "'0.1058"
$ws.Range("E42").Value = "'0.48%"
$ws.Range("E43").Value = "'-6.94%"
$ws.Range("D44").Value = "'0.009590"
$ws.Range("E44").Value = "'20.21%"
$ws.Range("D45").Value = "'0.00005104"
$ws.Range("E45").Value = "'-5.14%"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("D48").Value = "'0.002544"
$ws.Range("E48").Value = "'5.93%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.03%"
